$d = $word.ActiveDocument

# Locate the unique sentence fragment that contains the word to change,
# so we land on the correct occurrence of "Benutzer" (it also appears
# inside "Benutzeroberflaeche" earlier in the document).
$oldWord = "Benutzer"
$anchor = $d.Content
$anchor.Find.Execute("$oldWord ihre eigenen Analysen", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
$start = $anchor.Start

# Range covering just the word "Benutzer".
$target = $d.Range($start, $start + $oldWord.Length)

# Replace "Benutzer" with "Nutzende": delete the old word, then insert the
# new one. Toggling a character property on the inserted text and back
# off forces Word to keep it as its own run (matching how a manual
# type-over edit leaves the surrounding runs split apart) instead of
# silently re-merging it into the neighbouring runs.
$target.Delete()
$target.InsertAfter("Nutzende")
$target.Bold = 1
$target.Bold = 0
